# Update example data CON (DTI) - subject17 connectivity matrix
# Several asymmetric cell pairs (upper/lower triangle) are reconciled by
# overwriting one side of the pair with the value from the other side.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.82466700350286071
$ws.Range("C1").Value = 0.93506051133611856
$ws.Range("BO1").Value = 0.66392089361373463
$ws.Range("BP2").Value = 0.89335677986234918
$ws.Range("E3").Value = 0.67092188632537553
$ws.Range("BD3").Value = 0.81034321904513629
$ws.Range("B4").Value = 0.91644526624634559
$ws.Range("C4").Value = 0.90525867031213614
$ws.Range("D5").Value = 0.91446856174408186
$ws.Range("D6").Value = 0.92652031872168994
$ws.Range("E6").Value = 0.75058842388009983
$ws.Range("H6").Value = 0.71392281109811861
$ws.Range("E7").Value = 0.712049090792783
$ws.Range("F7").Value = 0.76441001018069343
$ws.Range("H7").Value = 0.84778050135621652
$ws.Range("AG7").Value = 0.92937894711528379
$ws.Range("I8").Value = 0.83119134872505729
$ws.Range("J8").Value = 0.81824899339809343
$ws.Range("AH8").Value = 0.7281042958195918
$ws.Range("G9").Value = 0.76286495613617578
$ws.Range("I11").Value = 0.93434059552160265
$ws.Range("K12").Value = 0.92932684858439185
$ws.Range("M12").Value = 0.93362675093466208
$ws.Range("N12").Value = 0.74399144134659445
$ws.Range("K13").Value = 0.99328928916569414
$ws.Range("N13").Value = 0.66342863385070805
$ws.Range("AC13").Value = 0.94134514584010565
$ws.Range("AX13").Value = 0.71070812536123262
$ws.Range("J14").Value = 0.92526852291478323
$ws.Range("O14").Value = 0.94059391970264594
$ws.Range("BJ15").Value = 0.98668698280987233
$ws.Range("Q16").Value = 0.95789808557353862
$ws.Range("R16").Value = 0.9021185330824214
$ws.Range("O17").Value = 0.89531136581611714
$ws.Range("S17").Value = 0.83526086773030439
$ws.Range("BL17").Value = 0.94854652313164989
$ws.Range("T18").Value = 0.94726771428971612
$ws.Range("BP18").Value = 0.92477338193199865
$ws.Range("U19").Value = 0.99718619061828107
$ws.Range("BP19").Value = 0.98318189263679667
$ws.Range("V20").Value = 0.73166794720449957
$ws.Range("AU20").Value = 0.90684345692749868
$ws.Range("BD21").Value = 0.6464116174890252
$ws.Range("U22").Value = 0.96099502767314016
$ws.Range("V23").Value = 0.72996917134424744
$ws.Range("AD23").Value = 0.97749316891685356
$ws.Range("BI23").Value = 0.78744278707254911
$ws.Range("V24").Value = 0.84756149768654243
$ws.Range("W25").Value = 0.63703092194881084
$ws.Range("AA25").Value = 0.98212814372237434
$ws.Range("E28").Value = 0.83186522738874369
$ws.Range("Z28").Value = 0.92997534482954247
$ws.Range("AQ28").Value = 0.7326153918557019
$ws.Range("BC28").Value = 0.71850743091106417
$ws.Range("AA29").Value = 0.7289352655505873
$ws.Range("AD29").Value = 0.95072350397032745
$ws.Range("AF30").Value = 0.79259556454103297
$ws.Range("AC31").Value = 0.84195108999252044
$ws.Range("AD31").Value = 0.9617311171177656
$ws.Range("Y32").Value = 0.85093971349125086
$ws.Range("AE32").Value = 0.84117335385113434
$ws.Range("AH32").Value = 0.8302198300356689
$ws.Range("AI32").Value = 0.98614215476483547
$ws.Range("AE33").Value = 0.79756911662041252
$ws.Range("AF33").Value = 0.86877071586947907
$ws.Range("AH33").Value = 0.89663323572364484
$ws.Range("AB35").Value = 0.71831231938243878
$ws.Range("AH35").Value = 0.91332471908605739
$ws.Range("AL36").Value = 0.88485117291585214
$ws.Range("I37").Value = 0.88751733915840414
$ws.Range("AJ37").Value = 0.89987575630154071
$ws.Range("AL37").Value = 0.69381538044754709
$ws.Range("AP37").Value = 0.62202160370207593
$ws.Range("BE38").Value = 0.63166421848769827
$ws.Range("J40").Value = 0.92316234019731258
$ws.Range("X40").Value = 0.71679814969095013
$ws.Range("AL40").Value = 0.8618917136500196
$ws.Range("AM40").Value = 0.61705987801289075
$ws.Range("AO40").Value = 0.94075664703578621
$ws.Range("K41").Value = 0.63604712286673426
$ws.Range("AM41").Value = 0.97106382247516632
$ws.Range("AQ41").Value = 0.91576533144714123
$ws.Range("AN42").Value = 0.57762674229765154
$ws.Range("AP43").Value = 0.94834942069752226
$ws.Range("BB43").Value = 0.99215393793519557
$ws.Range("BP43").Value = 0.86531850768672358
$ws.Range("AT44").Value = 0.775407234259057
$ws.Range("AQ45").Value = 0.93916912571615219
$ws.Range("AR45").Value = 0.96098495539965234
$ws.Range("AS46").Value = 0.94086790402846809
$ws.Range("AV46").Value = 0.8984675166378131
$ws.Range("AJ47").Value = 0.90682575267759646
$ws.Range("AS47").Value = 0.68802507517264522
$ws.Range("AT47").Value = 0.98607056747611121
$ws.Range("AW47").Value = 0.80207764297577278
$ws.Range("AX48").Value = 0.71600630981883717
$ws.Range("AV49").Value = 0.84857622405164723
$ws.Range("AX49").Value = 0.94858560961717875
$ws.Range("AY49").Value = 0.73976939327916935
$ws.Range("BC50").Value = 0.89823878207356489
$ws.Range("AM51").Value = 0.77142587875742985
$ws.Range("AY52").Value = 0.98916950224676659
$ws.Range("BA52").Value = 0.91426027713261959
$ws.Range("BB52").Value = 0.85651273762159086
$ws.Range("Q53").Value = 0.88311515607603486
$ws.Range("X54").Value = 0.88736652542451244
$ws.Range("AU54").Value = 0.88394971160240776
$ws.Range("BA54").Value = 0.80965223711856016
$ws.Range("T55").Value = 0.89630491529159029
$ws.Range("BA55").Value = 0.63858079611409002
$ws.Range("BF56").Value = 0.54143980502933542
$ws.Range("BD57").Value = 0.73436157221681375
$ws.Range("N58").Value = 0.93723334136093395
$ws.Range("BE58").Value = 0.92500025122001972
$ws.Range("BG58").Value = 0.77655315347601217
$ws.Range("X59").Value = 0.88492144691478869
$ws.Range("BE59").Value = 0.82155142390710423
$ws.Range("BI60").Value = 0.77678455984245964
$ws.Range("BJ60").Value = 0.60032835536416818
$ws.Range("BG61").Value = 0.66811295261333847
$ws.Range("BM61").Value = 0.92827104505839553
$ws.Range("BI62").Value = 0.65439971587443457
$ws.Range("A63").Value = 0.79631261692676392
$ws.Range("AF63").Value = 0.74224325402001212
$ws.Range("BJ63").Value = 0.67363466881383238
$ws.Range("B64").Value = 0.89738304929979451
$ws.Range("Z64").Value = 0.62393331906480198
$ws.Range("BJ64").Value = 0.51622098816946305
$ws.Range("BN64").Value = 0.71437321325248981
$ws.Range("AA66").Value = 0.98809367272803228
$ws.Range("BM66").Value = 0.82194912212470839
$ws.Range("BO66").Value = 0.99973303513788059
$ws.Range("BP66").Value = 0.57103300446304073
$ws.Range("AK67").Value = 0.85659835770909643
$ws.Range("BP67").Value = 0.96213047349495373
$ws.Range("A68").Value = 0.75363057107493181
